# Applies the "grammar etc" updates to practical_1.docx:
#  - four small wording tweaks in the learning-objectives / outline list
#  - rename the "instructions" bookmark to "excercise"
#  - retitle the "Instructions" Heading1 to "Excercise"

$d = $word.ActiveDocument

# 1. Learning objective bullet: "Models should be parsimonious..."
$d.Content.Find.Execute(
    "Models should be parsimonious, with as few parameters as possible to capture the dynamics of interest.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Be aware that models should be parsimonious, with as few parameters as possible to capture the dynamics of interest.",
    2) | Out-Null

# 2. Learning objective bullet: "Model structure may be subjective..."
$d.Content.Find.Execute(
    "Model structure may be subjective, there can be many approaches to answering a given question.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Understand that model structure may be subjective, there can be many approaches to answering a given question.",
    2) | Out-Null

# 3. Outline step: "Read through and choose a study question (5 minutes)"
$d.Content.Find.Execute(
    "Read through and choose a study question (5 minutes)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read through the excercise and choose a study question (5 minutes)",
    2) | Out-Null

# 4. Outline step: "Discuss possible for each study question (20 minutes)"
$d.Content.Find.Execute(
    "Discuss possible for each study question (20 minutes)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discuss possible model structures for each study question (20 minutes)",
    2) | Out-Null

# 5 & 6. Rename the "instructions" bookmark to "excercise" and retitle the
# heading it wraps from "Instructions" to "Excercise". Word has no direct
# bookmark-rename API, so re-create it under the new name at the same range
# and drop the old one.
$bm = $d.Bookmarks("instructions")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("excercise", $bmRange) | Out-Null

$d.Content.Find.Execute(
    "Instructions",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Excercise",
    2) | Out-Null
